$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("1").Delete()
$ws.Range("G7").Select()
